$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.27
$ws.Range("I3").Value = 10
$ws.Range("K3").Value = 2.4
$ws.Range("L3").Value = 10
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("Q3").Value = 1.85
$ws.Range("R3").Value = 2
$ws.Range("U3").Value = 2.5
$ws.Range("V3").Value = 1.5
$ws.Range("X3").Value = 5.5
$ws.Range("Y3").Value = 9.5
$ws.Range("AE3").Value = 29
$ws.Range("AF3").Value = 101
$ws.Range("AH3").Value = 19
$ws.Range("AJ3").Value = 29
$ws.Range("AL3").Value = 81
$ws.Range("AM3").Value = 81
$ws.Range("AO3").Value = 6
$ws.Range("AQ3").Value = 17
$ws.Range("AU3").Value = 11
$ws.Range("AZ3").Value = 301
$ws.Range("BA3").Value = 351

$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 3.4
$ws.Range("I6").Value = 1.85
$ws.Range("J6").Value = 4.33
$ws.Range("K6").Value = 2.1
$ws.Range("L6").Value = 2.5
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("W6").Value = 11
$ws.Range("X6").Value = 21
$ws.Range("Y6").Value = 15
$ws.Range("AA6").Value = 34
$ws.Range("AC6").Value = 9
$ws.Range("AE6").Value = 17
$ws.Range("AI6").Value = 8.5
$ws.Range("AJ6").Value = 9
$ws.Range("AK6").Value = 15
$ws.Range("AL6").Value = 17
$ws.Range("AN6").Value = 6
$ws.Range("AO6").Value = 23
$ws.Range("AP6").Value = 34
$ws.Range("AQ6").Value = 81
$ws.Range("AR6").Value = 101
$ws.Range("AW6").Value = 3.75
$ws.Range("AX6").Value = 10
$ws.Range("AY6").Value = 21
$ws.Range("AZ6").Value = 34

$wb.Save()
